# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows that changed after the data repull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -8
    5  = -7
    12 = -10
    17 = -5
    19 = -14
    20 = -13
    22 = -6
    27 = -9
    28 = 4
    42 = 0
    51 = 5
    55 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
